$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Structural edits: insert a new "Stage 0" column (B) and two new
#    parameter rows (traceMult, case_isolation) before the threshold rows.
# ---------------------------------------------------------------------------

# Insert a new column before B -- shifts old B:G -> C:H
$ws.Columns("B").Insert()

# Insert two new rows before row 12 -- shifts old rows 12:14 -> 14:16
$ws.Rows("12:13").Insert()

# Carry over the bordered cell-style (s="1") used throughout the table onto
# the freshly-inserted rows, reusing the existing style instead of minting a
# new one.
$ws.Range("A11:H11").Copy()
$ws.Range("A12:H13").PasteSpecial(-4122)

# The newly inserted column H (old G, which never had the bordered style)
# also needs the bordered style applied across the whole table.
$ws.Range("A1").Copy()
$ws.Range("H1:H16").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Stage 0"

# ---------------------------------------------------------------------------
# 3. New "Stage 0" column values (column B) for the existing rows
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = 20
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 100
$ws.Range("B7").Value = $true
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 25
$ws.Range("B10").Value = 0.1428
$ws.Range("B11").Value = 0.05

# ---------------------------------------------------------------------------
# 4. Value tweaks among the existing (now shifted) columns
# ---------------------------------------------------------------------------
$ws.Range("F7").Value = $true      # schoolsOpen, old column E (Stage 3->here Stage 2b) 0 -> 1

# ---------------------------------------------------------------------------
# 5. New parameter rows: traceMult (12) and case_isolation (13)
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "traceMult"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = "Multiplier for the proportion of people tracked each day, which itself varies based on case load."

$ws.Range("A13").Value = "case_isolation"
$ws.Range("B13").Value = $false
$ws.Range("C13").Value = $true
$ws.Range("D13").Value = $true
$ws.Range("E13").Value = $true
$ws.Range("F13").Value = $true
$ws.Range("G13").Value = $true
$ws.Range("H13").Value = "Whether tracked cases isolate and cause their household to isolate."

# ---------------------------------------------------------------------------
# 6. "NA" placeholders in the new Stage 0 column for the threshold rows
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = "NA"
$ws.Range("B15").Value = "NA"
$ws.Range("B16").Value = 2.5

# ---------------------------------------------------------------------------
# 7. Cosmetic bits: column widths, selection, window geometry
# ---------------------------------------------------------------------------
$ws.Range("B1").ColumnWidth = 24.140625
$ws.Range("H1").ColumnWidth = 94.7109375

$ws.Range("G25").Select()

$excel.ActiveWindow.WindowState = -4143
$wb.Windows.Item(1).Left = 30975
$wb.Windows.Item(1).Top = 945
$wb.Windows.Item(1).Width = 25035
$wb.Windows.Item(1).Height = 13980
